$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style from the last existing header cell (AE1) so the new headers match formatting
$ws.Range("AF1").Value = "%DiffH"
$ws.Range("AG1").Value = "%DiffD"
$ws.Range("AH1").Value = "%DiffA"

$ws.Range("AE1").Copy() | Out-Null
$ws.Range("AF1:AH1").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false
